{"js": "// Chapter 3.1.28 (\"N5. 48: Arkivdelreferanser\") gets an \"AND/OR\" alternative\n// answer appended after the existing Output paragraph, and the trailing\n// paragraph that only carried the stray _GoBack bookmark is dropped.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the existing \"Output\" answer paragraph by its text.\nlet outputPara = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ingen arkivdelreferanser er registrert\") !== -1) {\n    outputPara = items[i];\n    break;\n  }\n}\n\n// The trailing paragraph (holding only the _GoBack bookmark) is the last\n// paragraph in the body; drop it entirely.\nconst lastPara = items[items.length - 1];\nif (lastPara !== outputPara) {\n  lastPara.delete();\n  await context.sync();\n}\n\n// Append the new alternative-answer block right after the Output paragraph:\n//   (blank)\n//   AND/OR\n//   (blank)\n//   Arkivdelreferanser er registrert.\n//   (blank)\nlet cursor = outputPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"AND/OR\", Word.InsertLocation.after);\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"Arkivdelreferanser er registrert.\", Word.InsertLocation.after);\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Chapter 3.1.28 (\"N5. 48: Arkivdelreferanser\") gets an \"AND/OR\" alternative\n# answer appended after the existing Output paragraph, and the trailing\n# paragraph that only carried the stray _GoBack bookmark is dropped.\n\n$d = $word.ActiveDocument\n\n# Locate the existing \"Output\" answer paragraph by its text.\n$outputIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"Ingen arkivdelreferanser er registrert*\") {\n        $outputIndex = $i\n        break\n    }\n}\n\n# The trailing paragraph (holding only the _GoBack bookmark) is the last\n# paragraph in the body; drop it entirely.\n$lastIndex = $d.Paragraphs.Count\nif ($lastIndex -gt $outputIndex) {\n    $d.Paragraphs.Item($lastIndex).Range.Delete()\n}\n\n# Append the new alternative-answer block right after the Output paragraph:\n#   (blank)\n#   AND/OR\n#   (blank)\n#   Arkivdelreferanser er registrert.\n#   (blank)\n$anchor = $d.Paragraphs.Item($outputIndex).Range\n$anchor.InsertParagraphAfter()\n\n$p = $d.Paragraphs.Item($outputIndex + 1).Range\n$p.InsertParagraphAfter()\n\n$p2 = $d.Paragraphs.Item($outputIndex + 2).Range\n$p2.Text = \"AND/OR\"\n$p2.InsertParagraphAfter()\n\n$p3 = $d.Paragraphs.Item($outputIndex + 3).Range\n$p3.InsertParagraphAfter()\n\n$p4 = $d.Paragraphs.Item($outputIndex + 4).Range\n$p4.Text = \"Arkivdelreferanser er registrert.\"\n$p4.InsertParagraphAfter()\n"}
